$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column S header (link_referencia), matching header style of existing headers ---
$ws.Range("R1").Copy() | Out-Null
$ws.Range("S1").PasteSpecial(-4122) | Out-Null
$ws.Range("S1").Value = "link_referencia"
$excel.CutCopyMode = 0


# --- Row 2 ---
$ws.Range("I2").Value = "Depósito irregular de madeira"
$ws.Range("J2").Value = "A empresa mantinha em depósito 288,03 metros cúbicos de madeira em tora sem licença válida."
$c = $ws.Range("P2")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = $origStyle
$ws.Range("S2").Value = "http://jud-anexos.digesto.com.br/0cacd6d80c499ae25dcb85380a07c3dd.pdf"

# --- Row 3 ---
$ws.Range("I3").Value = "Desmatamento de vegetação nativa"
$ws.Range("J3").Value = "Destruição de 121,15 hectares de floresta nativa no bioma amazônico sem autorização da autoridade ambiental."
$ws.Range("S3").Value = "http://jud-anexos.digesto.com.br/90df1f9ac9917f6df4b9f91915b3a8bd.pdf"

# --- Row 4 ---
$ws.Range("I4").Value = "Queima de Rejeitos Industriais"
$ws.Range("Q4").Value = "NULL"
$ws.Range("S4").Value = "http://jud-anexos.digesto.com.br/209b98634c9c4a3f7c83d1c521c5b8d6.pdf"

# --- Row 5 ---
$ws.Range("J5").Value = "Ação civil pública movida por desmatamento e impedimento da regeneração da cobertura florestal explorada."
$ws.Range("P5").Value = 0
$ws.Range("S5").Value = "http://jud-anexos.digesto.com.br/7cef4bd654c356d84d7a617351f802cc.pdf"

# --- Row 6 ---
$ws.Range("J6").Value = "Supressão de castanheiras, poluição hídrica do Rio Itacaiúnas e Igarapé Salobo, impactos na Floresta Nacional do Tapirapé-Aquiri e nas comunidades indígenas Xikrin."
$ws.Range("L6").Value = "54000"
$ws.Range("M6").Value = "m2"
$ws.Range("N6").Value = $true
$ws.Range("O6").Value = "Obrigações de Fazer (com custo)"
$ws.Range("P6").Value = "0"
$ws.Range("S6").Value = "http://jud-anexos.digesto.com.br/d5d876cd4f1b1b52385f1ec4df9886ba.html"

# --- Row 7 ---
$ws.Range("S7").Value = "http://jud-anexos.digesto.com.br/24586eb78e174455c03d488d6518e16c.html"

# --- Row 8 ---
$ws.Range("P8").Value = 0
$ws.Range("S8").Value = "http://jud-anexos.digesto.com.br/be70db0dae0a1d1f18eb65dfd5523e76.html"

# --- Row 9 ---
$ws.Range("I9").Value = "Falha no fornecimento de energia elétrica"
$ws.Range("J9").Value = "Falhas no fornecimento de energia elétrica que causaram gastos com geradores e substituição de equipamentos."
$ws.Range("N9").Value = $true
$ws.Range("S9").Value = "http://jud-anexos.digesto.com.br/5631abc57e085f5121b0ead2d78e9e85.html"

# --- Row 10 ---
$ws.Range("I10").Value = "Ocupação Irregular de Terreno de Marinha"
$ws.Range("J10").Value = "Ocupação irregular de barraca de praia em terreno de marinha, sem licença, impactando o acesso público e o meio ambiente."
$ws.Range("M10").Value = "m2"
$ws.Range("P10").Value = 0
$ws.Range("S10").Value = "http://jud-anexos.digesto.com.br/a322a1e63179ad2a721b9268801dbfb6.pdf"

# --- Row 11 ---
$ws.Range("J11").Value = "Emissão de som em violação aos limites de horário e decibéis permitidos pela legislação, causando transtornos aos moradores do entorno."
$ws.Range("N11").Value = $false
$ws.Range("O11").Value = "NULL"
$c = $ws.Range("P11")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0"
$c.Style = $origStyle
$ws.Range("R11").Value = "NULL"
$ws.Range("S11").Value = "http://jud-anexos.digesto.com.br/f8941539197bdd583e3bdd448abb2e6c.html"

# --- Row 12 ---
$ws.Range("I12").Value = "Descumprimento de contrato"
$ws.Range("J12").Value = "Atraso na entrega de imóvel e divergências entre o que foi vendido e o que foi efetivamente entregue, causando prejuízos aos compradores."
$ws.Range("K12").Value = "NULL"
$ws.Range("N12").Value = $false
$ws.Range("S12").Value = "http://jud-anexos.digesto.com.br/c7fac8749343a02deddd5dfdcb5092e1.html"

# --- Row 13 ---
$ws.Range("I13").Value = "Danos ao sistema de saneamento"
$ws.Range("J13").Value = "Atraso na entrega de infraestrutura básica, como rede de água potável, resultando em uso de poço artesiano pelo comprador e falta de rede de esgoto sanitário."
$ws.Range("K13").Value = "07/03/2017"
$ws.Range("N13").Value = $false
$ws.Range("O13").Value = "NULL"
$ws.Range("S13").Value = "http://jud-anexos.digesto.com.br/179ce549bff6f6044d73c9c6f7b54d98.html"

# --- Row 14 ---
$ws.Range("I14").Value = "Abastecimento irregular de água"
$ws.Range("J14").Value = "Abastecimento irregular de água no bairro Jorge Texeira, com fornecimento de água de no máximo três horas diárias, passando a trinta minutos."
$ws.Range("O14").Value = "Custas Judiciais e Acordos"
$ws.Range("S14").Value = "http://jud-anexos.digesto.com.br/e50aea65864c6a6ca7ab8ab06025d54f.pdf"

# --- Row 15 ---
$ws.Range("I15").Value = "Abastecimento Irregular de Água"
$ws.Range("J15").Value = "Abastecimento de água irregular no bairro Jorge Teixeira, com fornecimento descontínuo e precário, afetando moradores e causando transtornos."
$ws.Range("S15").Value = "http://jud-anexos.digesto.com.br/c9c4529e46b75ee5e400e42350f2b2d7.pdf"

# --- Row 16 ---
$ws.Range("I16").Value = "Poluição Hídrica"
$ws.Range("J16").Value = "Interrupção no fornecimento de água por 3 dias devido ao rompimento de um encanamento."
$ws.Range("K16").Value = "20/02/2024"
$ws.Range("N16").Value = $false
$ws.Range("S16").Value = "http://jud-anexos.digesto.com.br/4a966dbc383fe11597026e3ca7432c93.pdf"

# --- Row 17 ---
$ws.Range("J17").Value = "Desmate de 11,826 hectares de Floresta Estacional Semidecidual, Bioma Mata Atlântica, sem autorização ambiental na Fazenda Manancial, zona rural de Águas Vermelhas."
$ws.Range("P17").Value = 0
$ws.Range("R17").Value = "NULL"
$ws.Range("S17").Value = "http://jud-anexos.digesto.com.br/e448073edc85a7c3dd335c89e2c64c2e.html"

# --- Row 18 ---
$ws.Range("I18").Value = "Desmatamento de vegetação nativa"
$ws.Range("J18").Value = "Supressão de 30 árvores nativas sem licença, depositadas em APP, impedindo regeneração. Regeneração natural ocorrida, mas compensação necessária."
$ws.Range("S18").Value = "http://jud-anexos.digesto.com.br/4af9d082c1a6dd5841a9c871eb76b26b.html"

# --- Row 19 ---
$ws.Range("J19").Value = "Captação irregular de água do Rio Verde para obras de loteamento, causando dano ambiental pontual e de pequena monta."
$ws.Range("S19").Value = "http://jud-anexos.digesto.com.br/cbdda0e55f226a6f1eac06b4a819d143.html"

# --- Row 20 ---
$ws.Range("I20").Value = "Colisão de veículo contra poste de energia elétrica"
$ws.Range("J20").Value = "O veículo colidiu contra um poste de energia elétrica na Praça Universitária, danificando o veículo locado."
$ws.Range("S20").Value = "http://jud-anexos.digesto.com.br/e748a98ec567a8446cf05d6974c9d6d5.html"

# --- Row 21 ---
$ws.Range("J21").Value = "Aterramento de nascente e intervenção em área de preservação permanente (APP) de aproximadamente 4 hectares e intervenção em curso d’água sem outorga."
$ws.Range("S21").Value = "http://jud-anexos.digesto.com.br/02597b35d5e8ece30d36627f30550386.html"

# --- Row 22 ---
$ws.Range("J22").Value = "Assoreamento do leito do Córrego Canta Galo devido a deslizamento de talude e erosão."
$c = $ws.Range("P22")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0"
$c.Style = $origStyle
$ws.Range("S22").Value = "http://jud-anexos.digesto.com.br/7d54d2fd28f74c625bb2157b3d507fe6.html"

# --- Row 23 ---
$ws.Range("S23").Value = "http://jud-anexos.digesto.com.br/fab3efa39a57f1d1a8bd1457980931ab.html"
